$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Nalco published one more day of price data (07-12-2025), which pushes every
# existing row in the table down by one. In this particular sheet, rows 2-35
# all still cite the same 02-11-2025 circular (price 296.05) so the net
# effect of "insert a new top row and shift everything down" is identical to
# "insert a duplicate row right after the last of those identical rows (at
# row 36) and advance the Date text in rows 2-35 by one calendar day" - that
# is the cheapest edit and is what we apply.

# Helper cell (well outside the used range) used to stage literal text so
# that writing dd-mm-yyyy-looking strings back into column A/E doesn't get
# auto-converted into real Excel dates: we format the helper as Text, copy
# it, and paste-special "values only" into the destination cell, which
# leaves the destination's own number format untouched.
$helper = $ws.Cells.Item(1000, 26)
$helper.NumberFormat = "@"

function Set-TextValue($cell, [string]$text) {
    $helper.Value2 = $text
    $helper.Copy()
    $cell.PasteSpecial(-4163) | Out-Null
}

# 1) Insert a new row at 36, duplicating row 35's content (same style as the
#    surrounding rows comes along automatically with Insert()).
$ws.Rows(36).Insert()

Set-TextValue $ws.Cells.Item(36, 1) "03-11-2025"
$ws.Cells.Item(36, 2).Value2 = "ALUMINIUM INGOT"
$ws.Cells.Item(36, 3).Value2 = "IE07"
$ws.Cells.Item(36, 4).Value2 = 296.05
Set-TextValue $ws.Cells.Item(36, 5) "02-11-2025"
$ws.Cells.Item(36, 6).Value2 = "https://nalcoindia.com/wp-content/uploads/2019/01/Ingot-02-11-2025.pdf"

# 2) Advance the "Date" column by one day for rows 2-35 (their other columns
#    are unchanged - same circular/price as before).
for ($r = 2; $r -le 35; $r++) {
    $cell = $ws.Cells.Item($r, 1)
    $old = $cell.Text
    $d = [DateTime]::ParseExact($old, "dd-MM-yyyy", $null)
    $new = $d.AddDays(1).ToString("dd-MM-yyyy")
    Set-TextValue $cell $new
}

$helper.Clear()
$excel.Application.CutCopyMode = $false
